$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A8: "in progress" -> "x"
$ws.Range("A8").Value = "x"

# D8: gets new text
$ws.Range("D8").Value = "alphabets only, no duplicates"

# New row 11
$ws.Range("C11").Value = "Restart game on failure"
$ws.Range("D11").Value = "When game ends restart the game when user"
$ws.Range("C11").VerticalAlignment = -4160
